$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.179.64"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "'1.865.69"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").Value = "'1.011"
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("D5").Value = "'313.10"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'1.008"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("D7").Value = "'0.5095"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("D8").Value = "'0.3882"
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "'0.08213"
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("D10").Value = "'1.112"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "'6.190"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "'1.851.47"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "'20.15"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").Value = "'7.190"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").Value = "'1.011"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "'90.57"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'0.06678"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "'17.60"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").Value = "'1.008"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "'5.982"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").Value = "'28.149.14"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").Value = "'2.208"
$ws.Range("E25").Value = "  -2.13%  "
$ws.Range("D26").Value = "'2.067.38"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").Value = "'158.93"
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("D28").Value = "'20.63"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("D29").Value = "'2.410"
$ws.Range("E29").Value = "  -3.95%  "
$ws.Range("D30").Value = "'125.69"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").Value = "'0.1042"
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("D32").Value = "'1.037"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("D34").Value = "'3.609"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "'9.277"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.06530"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02410"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").Value = "'0.6429"
$ws.Range("E39").Value = "  -2.53%  "
$ws.Range("D40").Value = "'1.243"
$ws.Range("E40").Value = "  +2.63%  "
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("D42").Value = "'4.944"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("D43").Value = "'11.15"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").Value = "'0.6042"
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("D45").Value = "'13.09"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").Value = "'3.669"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").Value = "'1.272"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("E48").Value = "  -2.03%  "
$ws.Range("D49").Value = "'1.204"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("D50").Value = "'121.01"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").Value = "'0.06887"
$ws.Range("E51").Value = "  +0.83%  "
